$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.550.12"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "2.003.58"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.92"
$ws.Range("E5").Value = "  -9.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.599"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.07"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.24"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0977"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.24"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.295.61"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.41"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.758"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.09"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "2.006.44"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "36.501.55"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.87"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "0.0₃0805"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.79"
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -8.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.16"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.67"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.128"
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.80"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.39"
$ws.Range("E33").Value = "  -4.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0607"
$ws.Range("E34").Value = "  -6.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.27"
$ws.Range("E36").Value = "  -5.70%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.69"
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.99"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0935"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.453.42"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0203"
$ws.Range("E44").Value = "  -4.61%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  -7.78%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.57"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.23"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.82"
$ws.Range("E49").Value = "  +24.15%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.86"
$ws.Range("E51").Value = "  -2.58%  "
